$d = $word.ActiveDocument

$d.Content.Find.Execute("71-52=19", $true, $false, $false, $false, $false, $true, 1, $false, "47+24=71", 2) | Out-Null
$d.Content.Find.Execute("87-62=25", $true, $false, $false, $false, $false, $true, 1, $false, "32+5=37", 2) | Out-Null
$d.Content.Find.Execute("25+20=45", $true, $false, $false, $false, $false, $true, 1, $false, "35+7=42", 2) | Out-Null
$d.Content.Find.Execute("2+69=71", $true, $false, $false, $false, $false, $true, 1, $false, "58-3=55", 2) | Out-Null
$d.Content.Find.Execute("86-38=48", $true, $false, $false, $false, $false, $true, 1, $false, "57+2=59", 2) | Out-Null
$d.Content.Find.Execute("82+3=85", $true, $false, $false, $false, $false, $true, 1, $false, "88-84=4", 2) | Out-Null
$d.Content.Find.Execute("42-21=21", $true, $false, $false, $false, $false, $true, 1, $false, "79-60=19", 2) | Out-Null
$d.Content.Find.Execute("12+68=80", $true, $false, $false, $false, $false, $true, 1, $false, "22+57=79", 2) | Out-Null
$d.Content.Find.Execute("1+76=77", $true, $false, $false, $false, $false, $true, 1, $false, "87-39=48", 2) | Out-Null
$d.Content.Find.Execute("35+0=35", $true, $false, $false, $false, $false, $true, 1, $false, "9+58=67", 2) | Out-Null
$d.Content.Find.Execute("36+51=87", $true, $false, $false, $false, $false, $true, 1, $false, "39-33=6", 2) | Out-Null
$d.Content.Find.Execute("6+88=94", $true, $false, $false, $false, $false, $true, 1, $false, "42-4=38", 2) | Out-Null
$d.Content.Find.Execute("40+51=91", $true, $false, $false, $false, $false, $true, 1, $false, "11+46=57", 2) | Out-Null
$d.Content.Find.Execute("83-47=36", $true, $false, $false, $false, $false, $true, 1, $false, "17+5=22", 2) | Out-Null
$d.Content.Find.Execute("28+25=53", $true, $false, $false, $false, $false, $true, 1, $false, "52+0=52", 2) | Out-Null
$d.Content.Find.Execute("71+21=92", $true, $false, $false, $false, $false, $true, 1, $false, "64-14=50", 2) | Out-Null
$d.Content.Find.Execute("83+16=99", $true, $false, $false, $false, $false, $true, 1, $false, "5+41=46", 2) | Out-Null
$d.Content.Find.Execute("73+7=80", $true, $false, $false, $false, $false, $true, 1, $false, "39+27=66", 2) | Out-Null
$d.Content.Find.Execute("85-74=11", $true, $false, $false, $false, $false, $true, 1, $false, "51+23=74", 2) | Out-Null
$d.Content.Find.Execute("54-35=19", $true, $false, $false, $false, $false, $true, 1, $false, "78-32=46", 2) | Out-Null
$d.Content.Find.Execute("92-60=32", $true, $false, $false, $false, $false, $true, 1, $false, "57+19=76", 2) | Out-Null
$d.Content.Find.Execute("98-24=74", $true, $false, $false, $false, $false, $true, 1, $false, "7+73=80", 2) | Out-Null
$d.Content.Find.Execute("44-15=29", $true, $false, $false, $false, $false, $true, 1, $false, "0+42=42", 2) | Out-Null
$d.Content.Find.Execute("97-16=81", $true, $false, $false, $false, $false, $true, 1, $false, "92-21=71", 2) | Out-Null
$d.Content.Find.Execute("70+3=73", $true, $false, $false, $false, $false, $true, 1, $false, "22+3=25", 2) | Out-Null
$d.Content.Find.Execute("35+63=98", $true, $false, $false, $false, $false, $true, 1, $false, "70-47=23", 2) | Out-Null
$d.Content.Find.Execute("45+50=95", $true, $false, $false, $false, $false, $true, 1, $false, "29+34=63", 2) | Out-Null
$d.Content.Find.Execute("59+1=60", $true, $false, $false, $false, $false, $true, 1, $false, "56-1=55", 2) | Out-Null
$d.Content.Find.Execute("6+37=43", $true, $false, $false, $false, $false, $true, 1, $false, "64-64=0", 2) | Out-Null
$d.Content.Find.Execute("49+44=93", $true, $false, $false, $false, $false, $true, 1, $false, "72+20=92", 2) | Out-Null
$d.Content.Find.Execute("19+77=96", $true, $false, $false, $false, $false, $true, 1, $false, "84-44=40", 2) | Out-Null
$d.Content.Find.Execute("45-22=23", $true, $false, $false, $false, $false, $true, 1, $false, "43-41=2", 2) | Out-Null
$d.Content.Find.Execute("4+68=72", $true, $false, $false, $false, $false, $true, 1, $false, "91-58=33", 2) | Out-Null
$d.Content.Find.Execute("59-32=27", $true, $false, $false, $false, $false, $true, 1, $false, "59-26=33", 2) | Out-Null
$d.Content.Find.Execute("17+30=47", $true, $false, $false, $false, $false, $true, 1, $false, "54-0=54", 2) | Out-Null
$d.Content.Find.Execute("23+0=23", $true, $false, $false, $false, $false, $true, 1, $false, "87-81=6", 2) | Out-Null
$d.Content.Find.Execute("21+9=30", $true, $false, $false, $false, $false, $true, 1, $false, "13+25=38", 2) | Out-Null
$d.Content.Find.Execute("8+88=96", $true, $false, $false, $false, $false, $true, 1, $false, "56+9=65", 2) | Out-Null
$d.Content.Find.Execute("73-23=50", $true, $false, $false, $false, $false, $true, 1, $false, "17+78=95", 2) | Out-Null
$d.Content.Find.Execute("24+43=67", $true, $false, $false, $false, $false, $true, 1, $false, "58-48=10", 2) | Out-Null
$d.Content.Find.Execute("86-10=76", $true, $false, $false, $false, $false, $true, 1, $false, "45+16=61", 2) | Out-Null
$d.Content.Find.Execute("95-46=49", $true, $false, $false, $false, $false, $true, 1, $false, "41-39=2", 2) | Out-Null
$d.Content.Find.Execute("7+4=11", $true, $false, $false, $false, $false, $true, 1, $false, "44-42=2", 2) | Out-Null
$d.Content.Find.Execute("37+15=52", $true, $false, $false, $false, $false, $true, 1, $false, "48+7=55", 2) | Out-Null
$d.Content.Find.Execute("17+10=27", $true, $false, $false, $false, $false, $true, 1, $false, "43-3=40", 2) | Out-Null
$d.Content.Find.Execute("95-14=81", $true, $false, $false, $false, $false, $true, 1, $false, "80-63=17", 2) | Out-Null
$d.Content.Find.Execute("81-3=78", $true, $false, $false, $false, $false, $true, 1, $false, "29+49=78", 2) | Out-Null
$d.Content.Find.Execute("46+35=81", $true, $false, $false, $false, $false, $true, 1, $false, "21-13=8", 2) | Out-Null
$d.Content.Find.Execute("44+29=73", $true, $false, $false, $false, $false, $true, 1, $false, "16+7=23", 2) | Out-Null
$d.Content.Find.Execute("2+75=77", $true, $false, $false, $false, $false, $true, 1, $false, "39+5=44", 2) | Out-Null
$d.Content.Find.Execute("56-30=26", $true, $false, $false, $false, $false, $true, 1, $false, "14+31=45", 2) | Out-Null
$d.Content.Find.Execute("97-81=16", $true, $false, $false, $false, $false, $true, 1, $false, "91-38=53", 2) | Out-Null
$d.Content.Find.Execute("13+7=20", $true, $false, $false, $false, $false, $true, 1, $false, "90-34=56", 2) | Out-Null
$d.Content.Find.Execute("23-16=7", $true, $false, $false, $false, $false, $true, 1, $false, "16+57=73", 2) | Out-Null
$d.Content.Find.Execute("45-9=36", $true, $false, $false, $false, $false, $true, 1, $false, "21+1=22", 2) | Out-Null
$d.Content.Find.Execute("14+44=58", $true, $false, $false, $false, $false, $true, 1, $false, "43+53=96", 2) | Out-Null
$d.Content.Find.Execute("8+41=49", $true, $false, $false, $false, $false, $true, 1, $false, "20+5=25", 2) | Out-Null
$d.Content.Find.Execute("32-16=16", $true, $false, $false, $false, $false, $true, 1, $false, "23-8=15", 2) | Out-Null
$d.Content.Find.Execute("45+37=82", $true, $false, $false, $false, $false, $true, 1, $false, "97-34=63", 2) | Out-Null
$d.Content.Find.Execute("13+13=26", $true, $false, $false, $false, $false, $true, 1, $false, "12+32=44", 2) | Out-Null
$d.Content.Find.Execute("96-78=18", $true, $false, $false, $false, $false, $true, 1, $false, "72+21=93", 2) | Out-Null
$d.Content.Find.Execute("90-38=52", $true, $false, $false, $false, $false, $true, 1, $false, "83-46=37", 2) | Out-Null
$d.Content.Find.Execute("83+10=93", $true, $false, $false, $false, $false, $true, 1, $false, "6+80=86", 2) | Out-Null
$d.Content.Find.Execute("68+23=91", $true, $false, $false, $false, $false, $true, 1, $false, "11+1=12", 2) | Out-Null
$d.Content.Find.Execute("47+3=50", $true, $false, $false, $false, $false, $true, 1, $false, "80-48=32", 2) | Out-Null
$d.Content.Find.Execute("77-61=16", $true, $false, $false, $false, $false, $true, 1, $false, "45+28=73", 2) | Out-Null
$d.Content.Find.Execute("21+76=97", $true, $false, $false, $false, $false, $true, 1, $false, "91+4=95", 2) | Out-Null
$d.Content.Find.Execute("99-57=42", $true, $false, $false, $false, $false, $true, 1, $false, "60-18=42", 2) | Out-Null
$d.Content.Find.Execute("89-18=71", $true, $false, $false, $false, $false, $true, 1, $false, "30+9=39", 2) | Out-Null
$d.Content.Find.Execute("24-11=13", $true, $false, $false, $false, $false, $true, 1, $false, "77-23=54", 2) | Out-Null
$d.Content.Find.Execute("98-29=69", $true, $false, $false, $false, $false, $true, 1, $false, "69-25=44", 2) | Out-Null
$d.Content.Find.Execute("17+46=63", $true, $false, $false, $false, $false, $true, 1, $false, "23+25=48", 2) | Out-Null
$d.Content.Find.Execute("52+2=54", $true, $false, $false, $false, $false, $true, 1, $false, "26+66=92", 2) | Out-Null
$d.Content.Find.Execute("54-8=46", $true, $false, $false, $false, $false, $true, 1, $false, "16+20=36", 2) | Out-Null
$d.Content.Find.Execute("49+28=77", $true, $false, $false, $false, $false, $true, 1, $false, "65-40=25", 2) | Out-Null
$d.Content.Find.Execute("93+1=94", $true, $false, $false, $false, $false, $true, 1, $false, "88-35=53", 2) | Out-Null
$d.Content.Find.Execute("39+32=71", $true, $false, $false, $false, $false, $true, 1, $false, "70+5=75", 2) | Out-Null
$d.Content.Find.Execute("82-56=26", $true, $false, $false, $false, $false, $true, 1, $false, "62-32=30", 2) | Out-Null
$d.Content.Find.Execute("8+19=27", $true, $false, $false, $false, $false, $true, 1, $false, "76-37=39", 2) | Out-Null
$d.Content.Find.Execute("83-4=79", $true, $false, $false, $false, $false, $true, 1, $false, "27+52=79", 2) | Out-Null
$d.Content.Find.Execute("14+60=74", $true, $false, $false, $false, $false, $true, 1, $false, "16-16=0", 2) | Out-Null
$d.Content.Find.Execute("25+49=74", $true, $false, $false, $false, $false, $true, 1, $false, "74-1=73", 2) | Out-Null
$d.Content.Find.Execute("24+26=50", $true, $false, $false, $false, $false, $true, 1, $false, "7+30=37", 2) | Out-Null
$d.Content.Find.Execute("85-57=28", $true, $false, $false, $false, $false, $true, 1, $false, "20+31=51", 2) | Out-Null
$d.Content.Find.Execute("27-19=8", $true, $false, $false, $false, $false, $true, 1, $false, "97-82=15", 2) | Out-Null
$d.Content.Find.Execute("28+49=77", $true, $false, $false, $false, $false, $true, 1, $false, "14+11=25", 2) | Out-Null
$d.Content.Find.Execute("63+27=90", $true, $false, $false, $false, $false, $true, 1, $false, "94-13=81", 2) | Out-Null
$d.Content.Find.Execute("1+71=72", $true, $false, $false, $false, $false, $true, 1, $false, "10+60=70", 2) | Out-Null
$d.Content.Find.Execute("62+34=96", $true, $false, $false, $false, $false, $true, 1, $false, "77+22=99", 2) | Out-Null
$d.Content.Find.Execute("75-6=69", $true, $false, $false, $false, $false, $true, 1, $false, "40-29=11", 2) | Out-Null
$d.Content.Find.Execute("11+62=73", $true, $false, $false, $false, $false, $true, 1, $false, "40+36=76", 2) | Out-Null
$d.Content.Find.Execute("41-34=7", $true, $false, $false, $false, $false, $true, 1, $false, "90+2=92", 2) | Out-Null
$d.Content.Find.Execute("96+3=99", $true, $false, $false, $false, $false, $true, 1, $false, "7+3=10", 2) | Out-Null
$d.Content.Find.Execute("68+24=92", $true, $false, $false, $false, $false, $true, 1, $false, "94-15=79", 2) | Out-Null
$d.Content.Find.Execute("98-62=36", $true, $false, $false, $false, $false, $true, 1, $false, "29+44=73", 2) | Out-Null
$d.Content.Find.Execute("88-83=5", $true, $false, $false, $false, $false, $true, 1, $false, "51+14=65", 2) | Out-Null
$d.Content.Find.Execute("11+57=68", $true, $false, $false, $false, $false, $true, 1, $false, "65-30=35", 2) | Out-Null
$d.Content.Find.Execute("41+40=81", $true, $false, $false, $false, $false, $true, 1, $false, "9+80=89", 2) | Out-Null
$d.Content.Find.Execute("20+67=87", $true, $false, $false, $false, $false, $true, 1, $false, "74+14=88", 2) | Out-Null
$d.Content.Find.Execute("9+5=14", $true, $false, $false, $false, $false, $true, 1, $false, "33+15=48", 2) | Out-Null
